$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cashOut (H) values that simply incremented by 1 ---
$ws.Cells.Item(3, 8).Value = 363.0
$ws.Cells.Item(5, 8).Value = 416.0
$ws.Cells.Item(7, 8).Value = 630.0
$ws.Cells.Item(15, 8).Value = 605.0

# --- Update cashOut (H) values that also pick up the "right aligned" style ---
# (these already matched an existing style once normalized, so re-apply the
# horizontal alignment explicitly to land on the de-duplicated style slot)
$ws.Cells.Item(9, 8).Value = 347.0
$ws.Cells.Item(9, 8).HorizontalAlignment = -4152

$ws.Cells.Item(11, 8).Value = 750.0
$ws.Cells.Item(11, 8).HorizontalAlignment = -4152

$ws.Cells.Item(13, 8).Value = 440.0
$ws.Cells.Item(13, 8).HorizontalAlignment = -4152

# --- Append three new poker-session rows (16-18), copying row 15's formats ---
$ws.Range("A15:H15").Copy()
$ws.Range("A16:H18").PasteSpecial(-4122)

$ws.Cells.Item(16, 1).Value = 43433.0
$ws.Cells.Item(16, 2).Value = "Thu"
$ws.Cells.Item(16, 3).Value = 2.0
$ws.Cells.Item(16, 4).Value = 43433.86388888889
$ws.Cells.Item(16, 5).Value = 43434.120833333334
$ws.Cells.Item(16, 6).Value = 500.0
$ws.Cells.Item(16, 7).Value = 0.0
$ws.Cells.Item(16, 8).Value = 693.0

$ws.Cells.Item(17, 1).Value = 43434.0
$ws.Cells.Item(17, 2).Value = "Fri"
$ws.Cells.Item(17, 3).Value = 3.0
$ws.Cells.Item(17, 4).Value = 43434.99930555555
$ws.Cells.Item(17, 5).Value = 43435.15625
$ws.Cells.Item(17, 6).Value = 300.0
$ws.Cells.Item(17, 7).Value = 0.0
$ws.Cells.Item(17, 8).Value = 724.0

$ws.Cells.Item(18, 1).Value = 43436.0
$ws.Cells.Item(18, 2).Value = "Mon"
$ws.Cells.Item(18, 3).Value = 2.0
$ws.Cells.Item(18, 4).Value = 43437.59375
$ws.Cells.Item(18, 5).Value = 43437.80902777778
$ws.Cells.Item(18, 6).Value = 300.0
$ws.Cells.Item(18, 7).Value = 0.0
$ws.Cells.Item(18, 8).Value = 851.0
